# Update E8 from "Good Morning" to "GIT UPDATE" and set the active cell / selection to E8.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

$ws.Activate()
$ws.Range("E8").Select()
